$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# Set the new cell value (zero string label) - numeric 0
$ws.Range("N1").Value = 0

# Add defined name for the new cell
$wb.Names.Add("zero_string", "='Lookups'!`$N`$1")
